# Fixed StudyComb for Faceted Filters ICDC
#
# The StatQuery (column C) used on the Cases/Samples/Files tabs is replaced
# with a corrected Cypher query. All three rows (2,3,4) share the exact same
# StatQuery text, so all three are updated identically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN [`"American Staffordshire Terrier`"]`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Column C width narrowed slightly and the saved view scroll/zoom settings
# were refreshed when the file was re-saved.
$ws.Columns.Item(3).ColumnWidth = 93.3

$ws.Application.ActiveWindow.Zoom = 100
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 1
